$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''317.59'
$ws.Range('E2').Value = '''3.44%'
$ws.Range('G2').Value = '''3'
$ws.Range('E3').Value = '''-0.37%'
$ws.Range('G3').Value = '''3'
$ws.Range('D4').Value = '''5.111'
$ws.Range('E4').Value = '''1.42%'
$ws.Range('G4').Value = '''3'
$ws.Range('D5').Value = '''0.08073'
$ws.Range('E5').Value = '''3.02%'
$ws.Range('G5').Value = '''3'
$ws.Range('D6').Value = '''2.159'
$ws.Range('E6').Value = '''1.80%'
$ws.Range('G6').Value = '''3'
$ws.Range('D7').Value = '''8.049'
$ws.Range('G7').Value = '''3'
$ws.Range('B8').Value = '''GateToken'
$ws.Range('C8').Value = '''https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range('D8').Value = '''4.129'
$ws.Range('E8').Value = '''1.68%'
$ws.Range('G8').Value = '''3'
$ws.Range('B9').Value = '''MXToken'
$ws.Range('C9').Value = '''https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D9').Value = '''0.9286'
$ws.Range('E9').Value = '''0.83%'
$ws.Range('G9').Value = '''3'
$ws.Range('B10').Value = '''LiechtensteinCryptoassetsExchange'
$ws.Range('C10').Value = '''https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range('D10').Value = '''0.1005'
$ws.Range('E10').Value = '''4.78%'
$ws.Range('G10').Value = '''3'
$ws.Range('B11').Value = '''WazirX'
$ws.Range('C11').Value = '''https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range('D11').Value = '''0.1884'
$ws.Range('E11').Value = '''0.04%'
$ws.Range('G11').Value = '''3'
$ws.Range('B12').Value = '''MandalaExchangeToken'
$ws.Range('C12').Value = '''https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range('D12').Value = '''0.09214'
$ws.Range('E12').Value = '''6.24%'
$ws.Range('G12').Value = '''3'
$ws.Range('B13').Value = '''BitrueCoin'
$ws.Range('C13').Value = '''https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range('D13').Value = '''0.03611'
$ws.Range('E13').Value = '''3.15%'
$ws.Range('G13').Value = '''3'
$ws.Range('B14').Value = '''BitMartToken'
$ws.Range('C14').Value = '''https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range('D14').Value = '''0.09941'
$ws.Range('E14').Value = '''0.10%'
$ws.Range('G14').Value = '''3'
$ws.Range('B15').Value = '''BitForexToken'
$ws.Range('C15').Value = '''https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range('D15').Value = '''0.001437'
$ws.Range('E15').Value = '''0.58%'
$ws.Range('G15').Value = '''3'
$ws.Range('B16').Value = '''TigerCash'
$ws.Range('C16').Value = '''https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range('D16').Value = '''0.005650'
$ws.Range('E16').Value = '''-1.05%'
$ws.Range('G16').Value = '''3'
$ws.Range('B17').Value = '''LEO'
$ws.Range('C17').Value = '''https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D17').Value = '''3.461'
$ws.Range('E17').Value = '''-0.01%'
$ws.Range('G17').Value = '''3'
$ws.Range('D18').Value = '''2.800'
$ws.Range('E18').Value = '''16.24%'
$ws.Range('G18').Value = '''3'
$ws.Range('E19').Value = '''-1.26%'
$ws.Range('G19').Value = '''3'
$ws.Range('D20').Value = '''0.1339'
$ws.Range('E20').Value = '''0.18%'
$ws.Range('G20').Value = '''3'
$ws.Range('D21').Value = '''5.070'
$ws.Range('E21').Value = '''6.25%'
$ws.Range('G21').Value = '''3'
$ws.Range('D22').Value = '''0.2202'
$ws.Range('G22').Value = '''3'
$ws.Range('D23').Value = '''0.04600'
$ws.Range('E23').Value = '''0.12%'
$ws.Range('G23').Value = '''3'
$ws.Range('D24').Value = '''0.001242'
$ws.Range('E24').Value = '''0.96%'
$ws.Range('G24').Value = '''3'
$ws.Range('D25').Value = '''0.004749'
$ws.Range('E25').Value = '''-6.90%'
$ws.Range('G25').Value = '''3'
$ws.Range('D26').Value = '''0.0001301'
$ws.Range('E26').Value = '''-6.98%'
$ws.Range('G26').Value = '''3'
$ws.Range('D27').Value = '''0.0004500'
$ws.Range('E27').Value = '''65.09%'
$ws.Range('G27').Value = '''3'
$ws.Range('G28').Value = '''3'
$ws.Range('G29').Value = '''3'
$ws.Range('G30').Value = '''3'
$ws.Range('G31').Value = '''3'
$ws.Range('G32').Value = '''3'
$ws.Range('G33').Value = '''3'
$ws.Range('G34').Value = '''3'
$ws.Range('G35').Value = '''3'
$ws.Range('G36').Value = '''3'
$ws.Range('G37').Value = '''3'
$ws.Range('G38').Value = '''3'
$ws.Range('D39').Value = '''0.01956'
$ws.Range('E39').Value = '''6.59%'
$ws.Range('G39').Value = '''3'
$ws.Range('D40').Value = '''0.04968'
$ws.Range('E40').Value = '''3.98%'
$ws.Range('G40').Value = '''3'
$ws.Range('D41').Value = '''0.007822'
$ws.Range('E41').Value = '''4.72%'
$ws.Range('G41').Value = '''3'
$ws.Range('D42').Value = '''0.1399'
$ws.Range('E42').Value = '''-0.28%'
$ws.Range('G42').Value = '''3'
$ws.Range('D43').Value = '''0.007689'
$ws.Range('E43').Value = '''-0.50%'
$ws.Range('G43').Value = '''3'
$ws.Range('D44').Value = '''0.002098'
$ws.Range('E44').Value = '''-6.18%'
$ws.Range('G44').Value = '''3'
$ws.Range('D45').Value = '''0.01166'
$ws.Range('E45').Value = '''12.25%'
$ws.Range('G45').Value = '''3'
$ws.Range('D46').Value = '''0.00006292'
$ws.Range('E46').Value = '''2.07%'
$ws.Range('G46').Value = '''3'
$ws.Range('D47').Value = '''0.00000000750'
$ws.Range('E47').Value = '''-0.04%'
$ws.Range('G47').Value = '''3'
$ws.Range('D48').Value = '''28.91'
$ws.Range('E48').Value = '''-32.42%'
$ws.Range('G48').Value = '''3'
$ws.Range('D49').Value = '''0.001900'
$ws.Range('E49').Value = '''-5.04%'
$ws.Range('G49').Value = '''3'
$ws.Range('D50').Value = '''0.00002100'
$ws.Range('E50').Value = '''-0.04%'
$ws.Range('G50').Value = '''3'
$ws.Range('D51').Value = '''0.0002000'
$ws.Range('E51').Value = '''-0.04%'
$ws.Range('G51').Value = '''3'
